$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.375.41"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.502.62"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'306.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").Value = "'96.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").Value = "'0.586"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "'36.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'7.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "'0.113"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.08%  "
$ws.Range("D14").Value = "2.882.52"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'15.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.39%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.501.13"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "'0.854"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "42.382.41"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "'12.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'6.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.94%  "
$ws.Range("D22").Value = "'71.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'253.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("D24").Value = "'2.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").Value = "'2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "'26.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +10.86%  "
$ws.Range("D29").Value = "'10.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'37.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.93%  "
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "'154.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'19.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.51%  "
$ws.Range("D34").Value = "'3.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "'0.0787"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("D36").Value = "'2.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  -5.53%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "'24.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.03%  "
$ws.Range("D41").Value = "'3.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'3.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'2.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'0.0301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "2.034.86"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "'84.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").Value = "2.741.50"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "'72.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.98%  "
$ws.Range("D51").Value = "'0.190"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
